$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$order = @("PTS", "REB", "AST", "STL", "BLK")

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
  $fCell = $ws.Cells.Item($row, 6)
  $gCell = $ws.Cells.Item($row, 7)

  $f = $fCell.Value2
  $g = $gCell.Value2

  # Reorder F: positions are PTS/REB/AST/BLK/STL -> PTS/REB/AST/STL/BLK (swap last two)
  $parts = $f.Split("/")
  $newF = $parts[0] + "/" + $parts[1] + "/" + $parts[2] + "/" + $parts[4] + "/" + $parts[3]

  # Reorder G: parse "LABEL val | LABEL val | ..." into map, then re-emit in new order
  $gparts = $g.Split(" | ")
  $map = @{}
  foreach ($p in $gparts) {
    $sp = $p.Split(" ")
    $map[$sp[0]] = $sp[1]
  }
  $outParts = @()
  foreach ($k in $order) {
    $outParts += "$k $($map[$k])"
  }
  $newG = $outParts -join " | "

  $fCell.Value = $newF
  $gCell.Value = $newG
}

Write-Output "done"
